# Auto-generated script applying the Ultros_Profits market-data refresh diff.
# For each changed cell we set the new numeric Value on the correct worksheet/cell.
# A handful of cells are fully removed by the refresh (no replacement value) -
# those are cleared explicitly instead of being given a numeric value.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ALC_cells = @(
    @{r=28; c=8; v=2095.0557},
    @{r=28; c=9; v=2141.1667},
    @{r=28; c=11; v=2141.1667},
    @{r=28; c=13; v=-1656.1667},
    @{r=33; c=8; v=476.69232},
    @{r=33; c=9; v=388.9565},
    @{r=33; c=10; v=1149.3334},
    @{r=33; c=11; v=388.9565},
    @{r=33; c=12; v=1149.3334},
    @{r=33; c=13; v=-159.9565},
    @{r=33; c=14; v=-1607.3334},
    @{r=75; c=8; v=55000},
    @{r=75; c=10; v=55000},
    @{r=75; c=12; v=55000},
    @{r=75; c=14; v=-56872},
    @{r=78; c=8; v=55000},
    @{r=78; c=10; v=55000},
    @{r=78; c=12; v=165000},
    @{r=78; c=14; v=-174360},
    @{r=106; c=8; v=5041.2144},
    @{r=106; c=9; v=5087.154},
    @{r=106; c=10; v=4444},
    @{r=106; c=11; v=5087.154},
    @{r=106; c=12; v=4444},
    @{r=106; c=13; v=-4456.154},
    @{r=106; c=14; v=-5706},
    @{r=112; c=8; v=1451.9181},
    @{r=112; c=10; v=1463.9656},
    @{r=112; c=12; v=4391.8968},
    @{r=112; c=14; v=-6607.8968},
    @{r=113; c=8; v=5927.722},
    @{r=113; c=9; v=3949.5},
    @{r=113; c=10; v=6175},
    @{r=113; c=11; v=3949.5},
    @{r=113; c=12; v=6175},
    @{r=113; c=13; v=-695.5},
    @{r=113; c=14; v=-12683},
    @{r=116; c=8; v=1998.75},
    @{r=116; c=9; v=1998.75},
    @{r=116; c=11; v=1998.75},
    @{r=116; c=13; v=1443.25},
    @{r=133; c=8; v=70000},
    @{r=133; c=10; v=70000},
    @{r=133; c=12; v=70000},
    @{r=133; c=14; v=-80120},
    @{r=137; c=8; v=2224.5715},
    @{r=137; c=9; v=1399.1428},
    @{r=137; c=10; v=3050},
    @{r=137; c=11; v=4197.428400000001},
    @{r=137; c=12; v=9150},
    @{r=137; c=13; v=-1647.428400000001},
    @{r=137; c=14; v=-14250},
    @{r=138; c=8; v=2720.8215},
    @{r=138; c=9; v=1712.9333},
    @{r=138; c=10; v=3883.7693},
    @{r=138; c=11; v=5138.7999},
    @{r=138; c=12; v=11651.3079},
    @{r=138; c=13; v=1.20010000000002},
    @{r=138; c=14; v=-21931.3079},
    @{r=141; c=8; v=3504.2},
    @{r=141; c=9; v=3318.6365},
    @{r=141; c=10; v=4865},
    @{r=141; c=11; v=9955.9095},
    @{r=141; c=12; v=14595},
    @{r=141; c=13; v=-4775.9095},
    @{r=141; c=14; v=-24955}
)
foreach ($cell in $ALC_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = $cell.v
}

$ws = $wb.Worksheets.Item("ARM")
$ARM_cells = @(
    @{r=2; c=8; v=13937.435},
    @{r=2; c=9; v=18126.705},
    @{r=2; c=11; v=18126.705},
    @{r=2; c=13; v=-18013.705},
    @{r=5; c=8; v=112.5},
    @{r=5; c=9; v=112.5},
    @{r=5; c=10; v=0},
    @{r=5; c=11; v=112.5},
    @{r=5; c=12; v=0},
    @{r=5; c=13; v=-0.5},
    @{r=12; c=8; v=1309.4},
    @{r=12; c=9; v=0},
    @{r=12; c=10; v=1309.4},
    @{r=12; c=11; v=0},
    @{r=12; c=12; v=1309.4},
    @{r=12; c=14; v=-1655.4},
    @{r=32; c=8; v=7501.242},
    @{r=32; c=9; v=5494.5537},
    @{r=32; c=11; v=5494.5537},
    @{r=32; c=13; v=-5207.5537},
    @{r=45; c=8; v=2425.9375},
    @{r=45; c=10; v=3037.75},
    @{r=45; c=12; v=3037.75},
    @{r=45; c=14; v=-3791.75},
    @{r=61; c=8; v=5223.154},
    @{r=61; c=9; v=2311.5},
    @{r=61; c=11; v=2311.5},
    @{r=61; c=13; v=-2099.5},
    @{r=107; c=8; v=75000},
    @{r=107; c=10; v=75000},
    @{r=107; c=12; v=75000},
    @{r=107; c=14; v=-82680},
    @{r=116; c=8; v=13937.435},
    @{r=116; c=9; v=18126.705},
    @{r=116; c=11; v=18126.705},
    @{r=116; c=13; v=-15832.705},
    @{r=136; c=8; v=5223.154},
    @{r=136; c=9; v=2311.5},
    @{r=136; c=11; v=6934.5},
    @{r=136; c=13; v=-4384.5}
)
foreach ($cell in $ARM_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = $cell.v
}
$ARM_clear_cells = @(
    @{r=5; c=14},
    @{r=12; c=13}
)
foreach ($cell in $ARM_clear_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = ""
}

$ws = $wb.Worksheets.Item("BSM")
$BSM_cells = @(
    @{r=3; c=8; v=13937.435},
    @{r=3; c=9; v=18126.705},
    @{r=3; c=11; v=18126.705},
    @{r=3; c=13; v=-18012.705},
    @{r=4; c=8; v=112.5},
    @{r=4; c=9; v=112.5},
    @{r=4; c=10; v=0},
    @{r=4; c=11; v=112.5},
    @{r=4; c=12; v=0},
    @{r=4; c=13; v=2.5},
    @{r=20; c=8; v=70039.8},
    @{r=20; c=9; v=2163.5},
    @{r=20; c=11; v=2163.5},
    @{r=20; c=13; v=-1916.5},
    @{r=94; c=8; v=2209.3809},
    @{r=94; c=9; v=2380.7778},
    @{r=94; c=10; v=1181},
    @{r=94; c=11; v=2380.7778},
    @{r=94; c=12; v=1181},
    @{r=94; c=13; v=-1929.7778},
    @{r=94; c=14; v=-2083},
    @{r=134; c=8; v=4316.393},
    @{r=134; c=9; v=3311.7727},
    @{r=134; c=11; v=9935.3181},
    @{r=134; c=13; v=-7400.3181}
)
foreach ($cell in $BSM_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = $cell.v
}
$BSM_clear_cells = @(
    @{r=4; c=14}
)
foreach ($cell in $BSM_clear_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = ""
}

$ws = $wb.Worksheets.Item("CRP")
$CRP_cells = @(
    @{r=16; c=8; v=525},
    @{r=16; c=9; v=387.66666},
    @{r=16; c=11; v=387.66666},
    @{r=16; c=13; v=-100.66666},
    @{r=22; c=8; v=250},
    @{r=22; c=9; v=250},
    @{r=22; c=11; v=250},
    @{r=22; c=13; v=100},
    @{r=60; c=8; v=5000},
    @{r=60; c=9; v=5000},
    @{r=60; c=11; v=5000},
    @{r=60; c=13; v=-4489},
    @{r=62; c=8; v=4891.7144},
    @{r=62; c=9; v=5548.7},
    @{r=62; c=10; v=3249.25},
    @{r=62; c=11; v=5548.7},
    @{r=62; c=12; v=3249.25},
    @{r=62; c=13; v=-4924.7},
    @{r=62; c=14; v=-4497.25},
    @{r=65; c=8; v=4891.7144},
    @{r=65; c=9; v=5548.7},
    @{r=65; c=10; v=3249.25},
    @{r=65; c=11; v=27743.5},
    @{r=65; c=12; v=16246.25},
    @{r=65; c=13; v=-24623.5},
    @{r=65; c=14; v=-22486.25},
    @{r=86; c=8; v=30993.908},
    @{r=86; c=9; v=41256},
    @{r=86; c=11; v=41256},
    @{r=86; c=13; v=-40133},
    @{r=89; c=8; v=30993.908},
    @{r=89; c=9; v=41256},
    @{r=89; c=11; v=206280},
    @{r=89; c=13; v=-200664},
    @{r=105; c=8; v=16673253},
    @{r=105; c=9; v=1887.5},
    @{r=105; c=11; v=1887.5},
    @{r=105; c=13; v=-140.5},
    @{r=113; c=8; v=525},
    @{r=113; c=9; v=387.66666},
    @{r=113; c=11; v=387.66666},
    @{r=113; c=13; v=1782.33334},
    @{r=132; c=8; v=2680.9583},
    @{r=132; c=9; v=2652},
    @{r=132; c=11; v=7956},
    @{r=132; c=13; v=-5426},
    @{r=141; c=8; v=106000},
    @{r=141; c=10; v=106000},
    @{r=141; c=12; v=106000},
    @{r=141; c=14; v=-116360}
)
foreach ($cell in $CRP_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = $cell.v
}

$ws = $wb.Worksheets.Item("CUL")
$CUL_cells = @(
    @{r=58; c=8; v=3810.4},
    @{r=58; c=10; v=4167.1113},
    @{r=58; c=12; v=12501.3339},
    @{r=58; c=14; v=-12757.3339},
    @{r=75; c=8; v=3786.6667},
    @{r=75; c=9; v=1445},
    @{r=75; c=10; v=4079.375},
    @{r=75; c=11; v=4335},
    @{r=75; c=12; v=12238.125},
    @{r=75; c=13; v=-3337},
    @{r=75; c=14; v=-14234.125},
    @{r=78; c=8; v=3786.6667},
    @{r=78; c=9; v=1445},
    @{r=78; c=10; v=4079.375},
    @{r=78; c=11; v=13005},
    @{r=78; c=12; v=36714.375},
    @{r=78; c=13; v=-8013},
    @{r=78; c=14; v=-46698.375},
    @{r=103; c=8; v=2445.6924},
    @{r=103; c=9; v=1895.1428},
    @{r=103; c=10; v=3088},
    @{r=103; c=11; v=5685.428400000001},
    @{r=103; c=12; v=9264},
    @{r=103; c=13; v=-4806.428400000001},
    @{r=103; c=14; v=-11022},
    @{r=110; c=8; v=8108.1665},
    @{r=110; c=9; v=6324.5},
    @{r=110; c=10; v=9000},
    @{r=110; c=11; v=18973.5},
    @{r=110; c=12; v=27000},
    @{r=110; c=13; v=-14883.5},
    @{r=110; c=14; v=-35180},
    @{r=132; c=8; v=1201.9546},
    @{r=132; c=9; v=1103.5},
    @{r=132; c=11; v=9931.5},
    @{r=132; c=13; v=-7401.5}
)
foreach ($cell in $CUL_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = $cell.v
}

$ws = $wb.Worksheets.Item("GSM")
$GSM_cells = @(
    @{r=70; c=8; v=89823.766},
    @{r=70; c=9; v=189052.83},
    @{r=70; c=10; v=4770.2856},
    @{r=70; c=11; v=189052.83},
    @{r=70; c=12; v=4770.2856},
    @{r=70; c=13; v=-188782.83},
    @{r=70; c=14; v=-5310.2856},
    @{r=73; c=8; v=89823.766},
    @{r=73; c=9; v=189052.83},
    @{r=73; c=10; v=4770.2856},
    @{r=73; c=11; v=189052.83},
    @{r=73; c=12; v=4770.2856},
    @{r=73; c=13; v=-188116.83},
    @{r=73; c=14; v=-6642.2856},
    @{r=122; c=8; v=4037.95},
    @{r=122; c=9; v=1590.5555},
    @{r=122; c=11; v=4771.666499999999},
    @{r=122; c=13; v=-2321.666499999999},
    @{r=138; c=8; v=69923},
    @{r=138; c=10; v=69923},
    @{r=138; c=12; v=69923},
    @{r=138; c=14; v=-80203}
)
foreach ($cell in $GSM_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = $cell.v
}

$ws = $wb.Worksheets.Item("LTW")
$LTW_cells = @(
    @{r=7; c=8; v=5535},
    @{r=7; c=9; v=6020.6875},
    @{r=7; c=10; v=4979.9287},
    @{r=7; c=11; v=6020.6875},
    @{r=7; c=12; v=4979.9287},
    @{r=7; c=13; v=-5908.6875},
    @{r=7; c=14; v=-5203.9287},
    @{r=46; c=8; v=438.92856},
    @{r=46; c=9; v=357.5},
    @{r=46; c=10; v=500},
    @{r=46; c=11; v=357.5},
    @{r=46; c=12; v=500},
    @{r=46; c=13; v=-169.5},
    @{r=46; c=14; v=-876},
    @{r=55; c=8; v=564.5625},
    @{r=55; c=9; v=219.28572},
    @{r=55; c=10; v=833.1111},
    @{r=55; c=11; v=219.28572},
    @{r=55; c=12; v=833.1111},
    @{r=55; c=13; v=-46.28572},
    @{r=55; c=14; v=-1179.1111},
    @{r=61; c=8; v=4479.6284},
    @{r=61; c=9; v=4457.8696},
    @{r=61; c=10; v=4521.3335},
    @{r=61; c=11; v=4457.8696},
    @{r=61; c=12; v=4521.3335},
    @{r=61; c=13; v=-4255.8696},
    @{r=61; c=14; v=-4925.3335},
    @{r=81; c=8; v=52666.668},
    @{r=81; c=10; v=52666.668},
    @{r=81; c=12; v=52666.668},
    @{r=81; c=14; v=-54662.668},
    @{r=84; c=8; v=52666.668},
    @{r=84; c=10; v=52666.668},
    @{r=84; c=12; v=158000.004},
    @{r=84; c=14; v=-167984.004},
    @{r=113; c=8; v=4479.6284},
    @{r=113; c=9; v=4457.8696},
    @{r=113; c=10; v=4521.3335},
    @{r=113; c=11; v=4457.8696},
    @{r=113; c=12; v=4521.3335},
    @{r=113; c=13; v=-2287.8696},
    @{r=113; c=14; v=-8861.333500000001},
    @{r=126; c=8; v=5535},
    @{r=126; c=9; v=6020.6875},
    @{r=126; c=10; v=4979.9287},
    @{r=126; c=11; v=18062.0625},
    @{r=126; c=12; v=14939.7861},
    @{r=126; c=13; v=-15592.0625},
    @{r=126; c=14; v=-19879.7861},
    @{r=136; c=8; v=5203.9546},
    @{r=136; c=9; v=3377.077},
    @{r=136; c=10; v=7842.778},
    @{r=136; c=11; v=10131.231},
    @{r=136; c=12; v=23528.334},
    @{r=136; c=13; v=-7581.231},
    @{r=136; c=14; v=-28628.334}
)
foreach ($cell in $LTW_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = $cell.v
}

$ws = $wb.Worksheets.Item("WVR")
$WVR_cells = @(
    @{r=16; c=8; v=30000},
    @{r=16; c=10; v=30000},
    @{r=16; c=12; v=30000},
    @{r=16; c=14; v=-30584},
    @{r=81; c=8; v=7611.1763},
    @{r=81; c=9; v=15906.286},
    @{r=81; c=11; v=31812.572},
    @{r=81; c=13; v=-30751.572},
    @{r=84; c=8; v=7611.1763},
    @{r=84; c=9; v=15906.286},
    @{r=84; c=11; v=159062.86},
    @{r=84; c=13; v=-153758.86},
    @{r=100; c=8; v=803.7},
    @{r=100; c=9; v=698.06665},
    @{r=100; c=11; v=1396.1333},
    @{r=100; c=13; v=-855.1333},
    @{r=107; c=8; v=445.33334},
    @{r=107; c=9; v=468.15384},
    @{r=107; c=11; v=1404.46152},
    @{r=107; c=13; v=515.5384799999999}
)
foreach ($cell in $WVR_cells) {
    $ws.Cells.Item($cell.r, $cell.c).Value = $cell.v
}
